$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All data rows (2 through 171) in column C ("Förändrad") were updated
# from serial date 45190 (2023-09-21) to 45192 (2023-09-23).
$ws.Range("C2:C171").Value = 45192
